$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-28: bump date serial value from 45184 to 45185
$ws.Range("C2:C28").Value = 45185
